$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("distance from Dma50")

$updates = @{
    "C2"  = 10.147
    "C3"  = 7.676
    "C4"  = 6.2524
    "C5"  = 5.2981
    "C6"  = 5.1282
    "C7"  = 4.9642
    "C8"  = 4.6458
    "C10" = 3.6659
    "C11" = 3.6647
    "C12" = 3.4922
    "C13" = 3.2317
    "C14" = 3.2288
    "C15" = 3.125
    "C16" = 3.0953
    "C17" = 2.8815
    "C18" = 2.6211
    "C19" = 2.5516
    "C20" = 2.3875
    "C21" = 2.3521
    "C22" = 1.4157
    "C23" = 1.3845
    "C24" = 1.2047
    "C25" = 1.1327
    "C26" = 0.9913999999999999
    "C27" = 0.9233
    "C28" = 0.6012999999999999
    "C29" = -0.1277
    "C30" = -2.059
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
